# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a "Label" column (H) to the worksheet, marking Control rows with 0
# and MDD rows with 1, and refreshes a handful of recomputed metric values
# that shifted slightly after the refit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" header in H1 ---
$ws.Range("H1").Value = "Label"

# --- Block 1: rows 2-11 (Iterations = 100) ---
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

# --- Block 2: rows 12-21 (Iterations = 200) ---
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1

# --- Refreshed metric values (refit results) ---
$ws.Range("D6").Value = 0.5617294745865449
$ws.Range("E6").Value = 0.5617294745865449

$ws.Range("D9").Value = 0.4697982634426429
$ws.Range("E9").Value = 0.5302017365573571

$ws.Range("F11").Value = 0.6968506574630737
